$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "230.85" into numbers,
# which would also drop formatting such as trailing zeros).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "43.990.83"
$ws.Range("E2").Value = "  +0.26%  "
Set-TextValue "D3" "2.266.73"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "230.85"
$ws.Range("E5").Value = "  -0.36%  "
Set-TextValue "D6" "0.630"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +6.08%  "
Set-TextValue "D10" "0.0998"
$ws.Range("E10").Value = "  +5.43%  "
Set-TextValue "D11" "57.33"
$ws.Range("E11").Value = "  -0.94%  "
Set-TextValue "D12" "27.28"
$ws.Range("E12").Value = "  +15.03%  "
$ws.Range("E13").Value = "  +1.78%  "
Set-TextValue "D14" "2.604.98"
$ws.Range("E14").Value = "  -0.51%  "
Set-TextValue "D15" "15.80"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("E16").Value = "  +4.95%  "
$ws.Range("E17").Value = "  +2.96%  "
Set-TextValue "D18" "2.261.75"
$ws.Range("E18").Value = "  -0.88%  "
Set-TextValue "D19" "43.888.22"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  +7.31%  "
Set-TextValue "D21" "73.80"
$ws.Range("E21").Value = "  +0.80%  "
Set-TextValue "D22" "6.11"
$ws.Range("E22").Value = "  -2.11%  "
Set-TextValue "D23" "252.60"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -4.10%  "
Set-TextValue "D26" "2.26"
$ws.Range("E26").Value = "  -4.52%  "
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "10.10"
$ws.Range("E27").Value = "  +2.14%  "
Set-TextValue "B28" "WEMIXToken"
Set-TextValue "C28" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D28" "3.26"
$ws.Range("E28").Value = "  +21.80%  "
Set-TextValue "D29" "171.37"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -1.36%  "
Set-TextValue "D31" "20.93"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E33").Value = "  +2.66%  "
Set-TextValue "D34" "0.0707"
$ws.Range("E34").Value = "  +6.87%  "
Set-TextValue "D35" "4.79"
$ws.Range("E35").Value = "  -0.34%  "
Set-TextValue "D36" "4.90"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("E37").Value = "  +5.23%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +0.15%  "
Set-TextValue "D43" "0.0991"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  +4.85%  "
Set-TextValue "D45" "8.25"
$ws.Range("E45").Value = "  -5.88%  "
Set-TextValue "D46" "10.42"
$ws.Range("E46").Value = "  +9.24%  "
$ws.Range("E47").Value = "  -0.72%  "
Set-TextValue "D48" "98.44"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -3.28%  "
Set-TextValue "B50" "NEARProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "2.33"
$ws.Range("E50").Value = "  +2.77%  "
Set-TextValue "B51" "Maker"
Set-TextValue "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D51" "1.447.34"
$ws.Range("E51").Value = "  -1.61%  "
